# Weekly CompStat report refresh (123rd Precinct) - "New crime data collected"
#
# Updates:
#  - Header: Volume/Number and the reporting week date range
#  - The weekly crime-stat grid (rows 16-30): new counts and recomputed
#    percent-change figures for the week, 28-day, YTD, 2yr/13yr/30yr columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: "Volume 30   Number  18" -> "...19"
#              "Report Covering the Week  5/1/2023  Through  5/7/2023"
#                -> "...5/8/2023  Through  5/14/2023"
# All runs inside these two strings share identical formatting, so replacing
# the whole text keeps the visual result identical to editing just the
# trailing run.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# ---------------------------------------------------------------------------
# Helper: write a numeric value into a cell, matching the number format of a
# reference cell that already carries the desired style (so no new cell
# style gets created).
# ---------------------------------------------------------------------------
function Set-NumFromRef($ws, $targetAddr, $refAddr, $value) {
    $ws.Range($targetAddr).NumberFormat = $ws.Range($refAddr).NumberFormat
    $ws.Range($targetAddr).Value = $value
}

# Helper: turn a cell back into one of the "placeholder" text cells (shared
# text "0" or "***.*") by copying an existing cell that already has exactly
# that look (value + style) onto the target.
function Copy-Cell($ws, $targetAddr, $refAddr) {
    $ws.Range($refAddr).Copy($ws.Range($targetAddr))
}

# ---------------------------------------------------------------------------
# Row 16 (Murder)
# ---------------------------------------------------------------------------
Set-NumFromRef $ws "C16" "F16" 1
Set-NumFromRef $ws "I16" "I16" 5
Set-NumFromRef $ws "K16" "K16" 25
Set-NumFromRef $ws "L16" "L16" 400
Set-NumFromRef $ws "M16" "M16" -37.5
Set-NumFromRef $ws "N16" "N16" -73.684210526315

# ---------------------------------------------------------------------------
# Row 17 (Rape)
# ---------------------------------------------------------------------------
Set-NumFromRef $ws "C17" "D17" 1
Set-NumFromRef $ws "E17" "E17" 0
Set-NumFromRef $ws "F17" "F17" 4
Set-NumFromRef $ws "H17" "H17" 100
Set-NumFromRef $ws "I17" "I17" 29
Set-NumFromRef $ws "J17" "J17" 15
Set-NumFromRef $ws "K17" "K17" 93.333333333333
Set-NumFromRef $ws "L17" "L17" 123.076923076923
Set-NumFromRef $ws "M17" "M17" 70.588235294117
Set-NumFromRef $ws "N17" "N17" -12.121212121212

# ---------------------------------------------------------------------------
# Row 18 (Robbery)
# ---------------------------------------------------------------------------
Set-NumFromRef $ws "F18" "F18" 2
Set-NumFromRef $ws "H18" "H18" -60
Set-NumFromRef $ws "I18" "I18" 21
Set-NumFromRef $ws "J18" "J18" 8
Set-NumFromRef $ws "K18" "K18" 162.5
Set-NumFromRef $ws "L18" "L18" 61.538461538461
Set-NumFromRef $ws "M18" "M18" -52.272727272727
Set-NumFromRef $ws "N18" "N18" -79.611650485436

# ---------------------------------------------------------------------------
# Row 19 (Fel. Assault)
# ---------------------------------------------------------------------------
Set-NumFromRef $ws "C19" "C19" 6
Set-NumFromRef $ws "D19" "D19" 8
Set-NumFromRef $ws "E19" "E19" -25
Set-NumFromRef $ws "G19" "G19" 19
Set-NumFromRef $ws "H19" "H19" 31.578947368421
Set-NumFromRef $ws "I19" "I19" 109
Set-NumFromRef $ws "J19" "J19" 100
Set-NumFromRef $ws "K19" "K19" 9
Set-NumFromRef $ws "L19" "L19" 81.666666666666
Set-NumFromRef $ws "M19" "M19" 105.660377358491
Set-NumFromRef $ws "N19" "N19" 81.666666666666

# ---------------------------------------------------------------------------
# Row 20 (Burglary) - C/D/E revert to the "no data" placeholder look
# ---------------------------------------------------------------------------
Copy-Cell $ws "C20" "C22"
Copy-Cell $ws "D20" "C22"
Copy-Cell $ws "E20" "E22"
Set-NumFromRef $ws "G20" "G20" 6
Set-NumFromRef $ws "H20" "H20" 16.666666666666
Set-NumFromRef $ws "I20" "I20" 29
Set-NumFromRef $ws "K20" "K20" -25.641025641025
Set-NumFromRef $ws "L20" "L20" 141.666666666667
Set-NumFromRef $ws "M20" "M20" 123.076923076923
Set-NumFromRef $ws "N20" "N20" -88.803088803088

# ---------------------------------------------------------------------------
# Row 21 (Gr. Larceny)
# ---------------------------------------------------------------------------
Set-NumFromRef $ws "C21" "C21" 8
Set-NumFromRef $ws "D21" "D21" 10
Set-NumFromRef $ws "E21" "E21" -20
Set-NumFromRef $ws "F21" "F21" 39
Set-NumFromRef $ws "H21" "H21" 18.181818181818
Set-NumFromRef $ws "I21" "I21" 194
Set-NumFromRef $ws "J21" "J21" 166
Set-NumFromRef $ws "K21" "K21" 16.867469879518
Set-NumFromRef $ws "L21" "L21" 94
Set-NumFromRef $ws "M21" "M21" 41.605839416058
Set-NumFromRef $ws "N21" "N21" -59.329140461215

# ---------------------------------------------------------------------------
# Row 24 (Transit)
# ---------------------------------------------------------------------------
Set-NumFromRef $ws "C24" "C24" 13
Set-NumFromRef $ws "D24" "D24" 11
Set-NumFromRef $ws "E24" "E24" 18.181818181818
Set-NumFromRef $ws "F24" "F24" 25
Set-NumFromRef $ws "G24" "G24" 40
Set-NumFromRef $ws "H24" "H24" -37.5
Set-NumFromRef $ws "I24" "I24" 177
Set-NumFromRef $ws "J24" "J24" 148
Set-NumFromRef $ws "K24" "K24" 19.594594594594
Set-NumFromRef $ws "L24" "L24" 113.253012048193
Set-NumFromRef $ws "M24" "M24" -5.851063829787

# ---------------------------------------------------------------------------
# Row 25 (Housing)
# ---------------------------------------------------------------------------
Set-NumFromRef $ws "D25" "D25" 1
Set-NumFromRef $ws "E25" "E25" 100
Set-NumFromRef $ws "G25" "G25" 12
Set-NumFromRef $ws "H25" "H25" -16.666666666666
Set-NumFromRef $ws "I25" "I25" 68
Set-NumFromRef $ws "K25" "K25" 0
Set-NumFromRef $ws "L25" "L25" 61.904761904761
Set-NumFromRef $ws "M25" "M25" -2.857142857142

# ---------------------------------------------------------------------------
# Row 27 (UCR Rape*)
# ---------------------------------------------------------------------------
Set-NumFromRef $ws "L27" "L27" -20

# ---------------------------------------------------------------------------
# Row 30 (Shooting Inc.) - D/E revert to the "no data" placeholder look
# ---------------------------------------------------------------------------
Copy-Cell $ws "D30" "C30"
Copy-Cell $ws "E30" "M30"
